$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "332.25"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "1.07%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "44.36"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "8.20%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.774"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "3.03%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08333"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "2.06%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "8.798"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.63%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "4.506"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-0.65%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.962"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-3.74%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.893"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-3.55%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9327"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "1.50%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1238"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-1.56%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.1945"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-0.53%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09424"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "1.01%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.03955"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "5.62%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.1066"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "1.06%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001309"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.58%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.005916"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-6.17%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.503"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "1.87%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "9.042"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "6.86%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1363"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-2.21%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2572"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "2.20%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04398"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.48%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001255"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-1.21%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004387"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "1.84%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001191"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "0.68%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0003994"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "0.04%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02836"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "3.48%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05703"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "5.44%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007931"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "3.48%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1425"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "1.04%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.009083"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-0.82%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002102"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-6.56%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.009949"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-11.75%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00007263"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "5.37%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000751"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.08%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003971"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "12.26%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.002281"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-0.15%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002103"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.08%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002003"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.08%"
